# estoy cambiando de lap
# Update "ciclo" (F column) values on the "student" sheet:
#   - rows 3:100   were "I"   -> now "II"
#   - rows 101:126 were "III" -> now "IV"
#   - rows 127:155 stay "VI" (unchanged)
# Also move the active selection from E157 to E13 (clears the stale
# topLeftCell/selection left over from the previous session).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("student")

$ws.Range("F3:F100").Value = "II"
$ws.Range("F101:F126").Value = "IV"

$ws.Range("E13").Select()
